$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "47.446.10"
$ws.Range("E2").Value = "  +5.55%  "

# Row 3
$ws.Range("D3").Value = "2.511.12"
$ws.Range("E3").Value = "  +3.38%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextValue "D5" "324.45"
$ws.Range("E5").Value = "  +2.16%  "

# Row 6
Set-TextValue "D6" "105.90"
$ws.Range("E6").Value = "  +1.82%  "

# Row 7
$ws.Range("E7").Value = "  +1.71%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
Set-TextValue "D9" "0.545"
$ws.Range("E9").Value = "  +2.83%  "

# Row 10
Set-TextValue "D10" "37.30"
$ws.Range("E10").Value = "  +4.34%  "

# Row 11
$ws.Range("E11").Value = "  +2.03%  "

# Row 12
Set-TextValue "D12" "0.123"
$ws.Range("E12").Value = "  +0.91%  "

# Row 13
Set-TextValue "D13" "18.52"
$ws.Range("E13").Value = "  +0.13%  "

# Row 14
$ws.Range("E14").Value = "  +4.15%  "

# Row 15
$ws.Range("D15").Value = "2.904.18"
$ws.Range("E15").Value = "  +3.55%  "

# Row 16
$ws.Range("D16").Value = "2.497.37"
$ws.Range("E16").Value = "  +3.36%  "

# Row 17
Set-TextValue "D17" "0.852"
$ws.Range("E17").Value = "  +2.41%  "

# Row 18
$ws.Range("D18").Value = "47.330.10"
$ws.Range("E18").Value = "  +5.73%  "

# Row 19
Set-TextValue "D19" "12.84"
$ws.Range("E19").Value = "  +3.83%  "

# Row 20
Set-TextValue "D20" "6.63"
$ws.Range("E20").Value = "  +4.14%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0945"
$ws.Range("E21").Value = "  +2.85%  "

# Row 22
Set-TextValue "D22" "71.11"
$ws.Range("E22").Value = "  +3.54%  "

# Row 23
Set-TextValue "D23" "253.20"
$ws.Range("E23").Value = "  +3.82%  "

# Row 24
$ws.Range("E24").Value = "  +4.18%  "

# Row 25
$ws.Range("E25").Value = "  +3.21%  "

# Row 26
Set-TextValue "D26" "26.51"
$ws.Range("E26").Value = "  +4.72%  "

# Row 27
$ws.Range("E27").Value = "  -0.10%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "2.30"
$ws.Range("E28").Value = "  +4.34%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D29" "10.09"
$ws.Range("E29").Value = "  +5.58%  "

# Row 30
Set-TextValue "D30" "35.47"
$ws.Range("E30").Value = "  +5.57%  "

# Row 31
Set-TextValue "D31" "0.134"
$ws.Range("E31").Value = "  +5.65%  "

# Row 32
Set-TextValue "D32" "49.72"
$ws.Range("E32").Value = "  +1.64%  "

# Row 33
Set-TextValue "D33" "19.82"
$ws.Range("E33").Value = "  +0.89%  "

# Row 34
Set-TextValue "D34" "5.33"
$ws.Range("E34").Value = "  +2.17%  "

# Row 35
Set-TextValue "D35" "0.0782"
$ws.Range("E35").Value = "  +2.19%  "

# Row 36
$ws.Range("E36").Value = "  +0.11%  "

# Row 37
$ws.Range("E37").Value = "  +2.55%  "

# Row 38
Set-TextValue "D38" "4.64"
$ws.Range("E38").Value = "  +4.22%  "

# Row 39
$ws.Range("E39").Value = "  +4.30%  "

# Row 40
Set-TextValue "D40" "123.36"
$ws.Range("E40").Value = "  -3.28%  "

# Row 41
$ws.Range("E41").Value = "  +2.05%  "

# Row 42
Set-TextValue "D42" "2.24"
$ws.Range("E42").Value = "  +3.13%  "

# Row 43
Set-TextValue "D43" "21.70"
$ws.Range("E43").Value = "  +3.39%  "

# Row 44
$ws.Range("E44").Value = "  +3.07%  "

# Row 45
$ws.Range("D45").Value = "1.983.66"
$ws.Range("E45").Value = "  +2.21%  "

# Row 46
Set-TextValue "D46" "3.04"
$ws.Range("E46").Value = "  +3.13%  "

# Row 47
$ws.Range("E47").Value = "  +0.66%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D48" "9.19"
$ws.Range("E48").Value = "  +0.20%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D49" "1.80"
$ws.Range("E49").Value = "  +1.32%  "

# Row 50
$ws.Range("E50").Value = "  +17.50%  "

# Row 51
Set-TextValue "D51" "79.80"
$ws.Range("E51").Value = "  +4.79%  "
